$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 31.824752
$ws.Range("H2").Value = 95.47425600000001
$ws.Range("I2").Value = 0.886907633630525
$ws.Range("J2").Value = 0.886907633630525
$ws.Range("M2").Value = 0.3045636666666667
$ws.Range("N2").Value = 0.913691
$ws.Range("Q2").Value = 9.692663159877334
$ws.Range("R2").Value = 87.23396843889601
$ws.Range("S2").Value = 0.886907633630525
$ws.Range("T2").Value = 0.886907633630525

# Row 3
$ws.Range("I3").Value = 0.06502043684278042
$ws.Range("J3").Value = 0.06502043684278042
$ws.Range("M3").Value = 0.3045636666666667
$ws.Range("N3").Value = 0.913691
$ws.Range("Q3").Value = 0.7105826682823334
$ws.Range("R3").Value = 6.395244014541001
$ws.Range("S3").Value = 0.06502043684278042
$ws.Range("T3").Value = 0.06502043684278042

# Row 4
$ws.Range("G4").Value = 1.696588
$ws.Range("H4").Value = 5.089764000000001
$ws.Range("I4").Value = 0.04728133775640876
$ws.Range("J4").Value = 0.04728133775640876
$ws.Range("M4").Value = 0.3045636666666667
$ws.Range("N4").Value = 0.913691
$ws.Range("Q4").Value = 0.5167190621026667
$ws.Range("R4").Value = 4.650471558924001
$ws.Range("S4").Value = 0.04728133775640876
$ws.Range("T4").Value = 0.04728133775640876

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.02836866666666667
$ws.Range("H5").Value = 0.085106
$ws.Range("I5").Value = 0.0007905917702857979
$ws.Range("J5").Value = 0.0007905917702857978
$ws.Range("M5").Value = 0.3045636666666667
$ws.Range("N5").Value = 0.913691
$ws.Range("Q5").Value = 0.008640065138444444
$ws.Range("R5").Value = 0.077760586246
$ws.Range("S5").Value = 0.0007905917702857979
$ws.Range("T5").Value = 0.0007905917702857978
